# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition listing) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 41
$ws1.Range("F8").Value  = 13
$ws1.Range("F9").Value  = 8196
$ws1.Range("F10").Value = 769
$ws1.Range("F11").Value = 259
$ws1.Range("F12").Value = 1110
$ws1.Range("F13").Value = 817
$ws1.Range("F14").Value = 45
$ws1.Range("F15").Value = 35
$ws1.Range("F16").Value = 207
$ws1.Range("F17").Value = 92
$ws1.Range("F19").Value = 212
$ws1.Range("F20").Value = 880

# --- Sheet "全部类型" (all categories listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 41
$ws4.Range("F10").Value = 13
$ws4.Range("F11").Value = 8196
$ws4.Range("F12").Value = 769
$ws4.Range("F13").Value = 259
$ws4.Range("F14").Value = 1110
$ws4.Range("F15").Value = 817
$ws4.Range("F16").Value = 45
$ws4.Range("F17").Value = 35
$ws4.Range("F18").Value = 207
$ws4.Range("F19").Value = 92
$ws4.Range("F21").Value = 212
$ws4.Range("F22").Value = 880
